$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: BandPassFilterKa
$ws.Range("A24").Value = "BandPassFilterKa"

# Row 25: description text (entered before the row 24 hyperlink URL so that
# the shared-string table matches the order the cells were authored in)
$ws.Range("E25").Value = "1.3dB IL"

# Row 24: hyperlink in column E
$ws.Range("E24").Value = "https://www.minicircuits.com/WebStore/dashboard.html?model=ZVBP-40600-K%2B"
$ws.Hyperlinks.Add($ws.Range("E24"), "https://www.minicircuits.com/WebStore/dashboard.html?model=ZVBP-40600-K%2B")
$ws.Range("E24").Style = "Hyperlink"

# Row 26: HPA KA
$ws.Range("A26").Value = "HPA KA"
$ws.Range("B26").Value = 35
$ws.Range("C26").Value = 7
$ws.Range("E26").Value = "https://rfessentials.com/full-ka-band-high-power-amplifier-26-ghz-to-40-ghz/"

# Row 27: Mixer
$ws.Range("A27").Value = "Mixer"
$ws.Range("B27").Value = -8
$ws.Range("C27").Value = 8
$ws.Range("E27").Value = "https://www.qorvo.com/products/p/CMD313"

# Resize column A to fit the new content and update the selection to reflect
# where the user clicked next.
$ws.Columns("A").AutoFit()
$ws.Range("A28").Select()
